$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The latest automated sync nudged the stored serial for the previous
# reading's timestamp by a hair of floating-point precision.
$ws.Cells.Item(24, 1).Value = 45876.91711516204

# Append the new sensor reading as row 25.
$ws.Cells.Item(25, 1).Value = 45876.95855033756
$ws.Cells.Item(25, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(25, 2).Value = 2025
$ws.Cells.Item(25, 3).Value = 28
$ws.Cells.Item(25, 4).Value = 14.22
$ws.Cells.Item(25, 5).Value = 91.25
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 2.94
$ws.Cells.Item(25, 8).Value = "ESE"
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = "23:00:18"
